# Update column C (Fitness) values on Sheet1 according to the commit diff.
# Current values are all 7310; new values vary by row-range as below.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Each entry: start row, end row, new value (applies to column C)
$ranges = @(
    @(2,   9,   8594),
    @(10,  10,  8523),
    @(11,  11,  8402),
    @(12,  28,  8176),
    @(29,  33,  8175),
    @(34,  55,  7836),
    @(56,  66,  7828),
    @(67,  116, 7639),
    @(117, 192, 7320),
    @(193, 202, 7312),
    @(203, 252, 7293)
)

foreach ($r in $ranges) {
    $startRow = $r[0]
    $endRow   = $r[1]
    $value    = $r[2]
    $rangeAddr = "C$startRow`:C$endRow"
    $ws.Range($rangeAddr).Value = $value
}
